# BankRmtf.xlsx maintenance edit
#
# Content changes being applied (per the commit's XML diff):
#  1. Sheet "DBD": the "TitaTlrNo" (row 26, 經辦) and "TitaTxtNo"
#     (row 27, 交易序號) field rows get a new remark in column G:
#     "入帳後更新" (was blank before).
#  2. Sheet "DBS": a new lookup-helper row is appended describing the
#     "findEntryDate" function and its condition text
#     "EntryDate >= ,AND AcDate <= ".
#
# (The remainder of the published XML diff — xr:* revision GUIDs,
# fileVersion/rupBuild bump, extra theme fonts, workbookView/selection
# state, de-duplicated cellXfs entries, absPath drive letter — are all
# artifacts of re-saving the workbook with a newer Excel build and are
# reproduced automatically by the host application; they carry no
# spreadsheet content of their own.)

$wb = $excel.ActiveWorkbook

# --- Sheet "DBD": add "入帳後更新" remarks -----------------------------
$dbd = $wb.Worksheets.Item("DBD")

$dbd.Range("G26").Value = "入帳後更新"
$dbd.Range("G27").Value = "入帳後更新"

# --- Sheet "DBS": add findEntryDate helper row -------------------------
$dbs = $wb.Worksheets.Item("DBS")

$dbs.Range("A3").Value = "findEntryDate"
$dbs.Range("B3").Value = "EntryDate >= ,AND AcDate <= "
